# Append: 2025-10-04 06:24 JST
# Update the top two rows of the "ランサーズ" sheet with the freshly
# scraped entries and drop the rest of the previous scrape's rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 2: new top entry ---------------------------------------------
$ws.Range("A2").Value = "2025-10-04 06:24:01"
$ws.Range("B2").Value = "【せどり×ツール製作】APIを使用したせどりツールを製作できるエンジニアさんを募集します♪"
$ws.Range("D2").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("G2").Value = 243
$ws.Range("H2").Value = "🔥API ◆ツール"

# --- Row 3: new second entry -------------------------------------------
$ws.Range("A3").Value = "2025-10-04 06:24:01"
$ws.Range("B3").Value = "【個人利用】Web情報収集の仕組み構築ご依頼"
$ws.Range("D3").Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Range("G3").Value = 10
$ws.Range("H3").ClearContents()

# --- Drop the old rows 4-21 (only 2 entries remain now) -----------------
$ws.Range("A4:H21").EntireRow.Delete()

# --- Hyperlinks: rebuild for just F2/F3 with the fresh URLs -------------
$ws.Range("F2:F3").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5217096")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5406440")
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5217096"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5406440"
$ws.Range("F2").Style = "Hyperlink"
$ws.Range("F3").Style = "Hyperlink"

# --- Column width tweaks --------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 46.166666666666664
$ws.Columns.Item(4).ColumnWidth = 25.166666666666668
$ws.Columns.Item(8).ColumnWidth = 11.166666666666666
